$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells: "<name>_old" -> "<name>_FV2310", "<name>_new" -> "<name>_FV2404" ---
$headerBases = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $headerBases.Count; $i++) {
    # Columns A..J (1..10) carry the "_old" -> "_FV2310" headers
    $ws.Cells.Item(1, $i + 1).Value = "$($headerBases[$i])_FV2310"
    # Columns L..U (12..21) carry the "_new" -> "_FV2404" headers (column K = "diff", unchanged)
    $ws.Cells.Item(1, $i + 12).Value = "$($headerBases[$i])_FV2404"
}

# --- 2. Turn the used range into an Excel Table ("Table1") with autofilter + banded rows ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U64"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3. Freeze the header row (top row) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "edit complete"
